$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before the existing "District" column (F),
# shifting the old F (District) data to G.
$ws.Range("F1").EntireColumn.Insert()

# New column header
$ws.Range("F2").Value = "Address"

# Fill in the new "Address" column values (derived from the school name
# portion of column B, the second line of the combined name+address text).
$ws.Range("F3").Value = "Govt. High School Chincholi(H) Chittapur"
$ws.Range("F4").Value = "Navodaya P U College Channarayapatna"
$ws.Range("F5").Value = "G H S NalwarChittapur"
$ws.Range("F6").Value = "Govt. High School GundgurthiChittapur"
$ws.Range("F7").Value = "Shri Jagadguru Gangadhar High SchoolHubballi"
$ws.Range("F8").Value = "S B H S HulikalMagadi"
$ws.Range("F9").Value = "L N J P N H S AlwaiBhalki"
$ws.Range("F10").Value = "St. Josephs High SchoolKoppa"
$ws.Range("F11").Value = "Sri Sadguru High School BasarikatteKoppa"
$ws.Range("F12").Value = "Shivaji H S Bhalki"
$ws.Range("F13").Value = "G J C ShravanabelagolaC R Patna"
$ws.Range("F14").Value = "Fatima High School KeshwapurHubli"
$ws.Range("F15").Value = "N E H S GalagaliBilagi"
$ws.Range("F16").Value = "Vedhavathi Girls High School Kadur"
$ws.Range("F17").Value = "G J C SingrihalliHarapanahalli"
$ws.Range("F18").Value = "S S R High School KanchikereHarapanahalli"
$ws.Range("F19").Value = "G H S Wadi(Jn)Chittapur"
$ws.Range("F20").Value = "G H S HanumidiBelur"
$ws.Range("F21").Value = "G J C Harapanahalli"
$ws.Range("F22").Value = "Fatima High School Hubballi"
$ws.Range("F23").Value = "Girls English School Hubballi"
$ws.Range("F24").Value = "L B H S HirekodigeKoppa"
$ws.Range("F25").Value = "Sri Siddaganga Rural High School GudemaranahalliMagadi"
$ws.Range("F26").Value = "G H S KundurmuttC R Patna"
$ws.Range("F27").Value = "G H S IngalagiChittapur"
$ws.Range("F28").Value = "S G V V D S High School ArasikereHarapanahalli"
$ws.Range("F29").Value = "National High School ChilurHonnali"
$ws.Range("F30").Value = "G H S Gulasindha Channarayapatna"
$ws.Range("F31").Value = "Govt. High School HebbalChittapur"
$ws.Range("F32").Value = "J H S HiremegalagereHarapanahalli"
$ws.Range("F33").Value = "S S P S High School NarveKoppa"
$ws.Range("F35").Value = "Chittapur"
$ws.Range("F36").Value = "Savarad S R H S TumbagiTalikoti"
$ws.Range("F37").Value = "Govt. Boy’s High School ShahabadChittapur"
$ws.Range("F38").Value = "Dr. Raj Kumar High School Shikaripura"
$ws.Range("F39").Value = "G E S Deshpandenagar Hubli"
$ws.Range("F40").Value = "Govt. High School Kudregundi"
$ws.Range("F41").Value = "Kalburgi South"
$ws.Range("F42").Value = "G H S Halahalli(K)Bhalki"
$ws.Range("F43").Value = "S S R H S ByalakereMagadi"
$ws.Range("F44").Value = "Govt. P U College (High School Section)B H Road"
$ws.Range("F45").Value = "Jayapura High SchoolKoppa"
$ws.Range("F46").Value = "G H S K Byrapura C R Patna"
$ws.Range("F47").Value = "S S R H S K B Math Magadi"
$ws.Range("F48").Value = "Govt. High School RummangudChincholi"
$ws.Range("F49").Value = "G H S LakhangaonBhalki"
$ws.Range("F50").Value = "G H S NiluvagiluKoppa"
$ws.Range("F51").Value = "Sri Allamaprabhu High School Balligari"
$ws.Range("F53").Value = "Sri Maruthi High School HolalurShivamogga"
$ws.Range("F54").Value = "R H C H S BagaliHarapanahalli"
$ws.Range("F55").Value = "G H S HalbargaBhalki"
$ws.Range("F56").Value = "G H S RawoorChittapur"
$ws.Range("F57").Value = "A P H S TheligiHarapanahalli"
$ws.Range("F58").Value = "V B C High School Muddebihal"
$ws.Range("F59").Value = "K L E S C P High School MahalingapurMudhol"
